$d = $word.ActiveDocument

# --- First paragraph: add a paragraph border (space-only, no line) ---
$para1 = $d.Paragraphs(1)
$para1.Range.Borders.DistanceFromTop = 5
$para1.Range.Borders.DistanceFromLeft = 5
$para1.Range.Borders.DistanceFromBottom = 5
$para1.Range.Borders.DistanceFromRight = 5

# --- First paragraph: change left indent from 120 twips (6pt) to 225 twips (11.25pt) ---
$para1.Format.LeftIndent = 11.25

# --- First paragraph text: drop the trailing " " run, then rename the placeholder ID ---
$spaceRange = $d.Range(31, 32)
$spaceRange.Delete()

$idRange = $d.Range(0, 31)
$null = $idRange.Find.Execute("**ID__AFFARS_5323_topic_2__ID**", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "**ID__AFFARS_SUBPART_5323_3__ID**", 2)
